$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MD SPA_Last counter (column B, row 2) from 9362 to 9371
$ws.Range("B2").Value = 9371

# Move the active selection to E2 to match the final saved state
$ws.Range("E2").Select()
